$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$ws.Range("A5").Value = "2026-02-09 05:33:41"
$ws.Range("B5").Value = "Ibrahim Usman Umar "

# "22" looks numeric, but the source data stores it as text - use a leading
# apostrophe to force text entry, then clear the resulting style override so
# the cell keeps the same (default) style as its neighbours.
$ws.Range("C5").Value = "'22"
$ws.Range("C5").Style = "Normal"

$ws.Range("D5").Value = 7
